$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 422.1579
$ws.Range("I19").Value = 276
$ws.Range("J19").Value = 623.125
$ws.Range("K19").Value = 276
$ws.Range("L19").Value = 623.125
$ws.Range("M19").Value = -101
$ws.Range("N19").Value = -973.125

$ws.Range("H116").Value = 2582.3794
$ws.Range("I116").Value = 2563.9167
$ws.Range("J116").Value = 2671
$ws.Range("K116").Value = 2563.9167
$ws.Range("L116").Value = 2671
$ws.Range("M116").Value = 878.0832999999998
$ws.Range("N116").Value = -9555

$ws.Range("H132").Value = 16396794
$ws.Range("I132").Value = 22730676
$ws.Range("J132").Value = 3217.2942
$ws.Range("K132").Value = 68192028
$ws.Range("L132").Value = 9651.882599999999
$ws.Range("M132").Value = -68189498
$ws.Range("N132").Value = -14711.8826

$ws.Range("H137").Value = 1712.6333
$ws.Range("I137").Value = 1331.3846
$ws.Range("J137").Value = 2004.1765
$ws.Range("K137").Value = 3994.1538
$ws.Range("L137").Value = 6012.529500000001
$ws.Range("M137").Value = -1444.1538
$ws.Range("N137").Value = -11112.5295

$ws.Range("H138").Value = 597088.0600000001
$ws.Range("I138").Value = 1000
$ws.Range("J138").Value = 634343.5600000001
$ws.Range("K138").Value = 3000
$ws.Range("L138").Value = 1903030.68
$ws.Range("M138").Value = 2140
$ws.Range("N138").Value = -1913310.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""  # remove cell (was -20518)

$ws.Range("H32").Value = 2938.4924
$ws.Range("I32").Value = 2920.0793
$ws.Range("K32").Value = 2920.0793
$ws.Range("M32").Value = -2633.0793

$ws.Range("H68").Value = 46666
$ws.Range("I68").Value = 30000
$ws.Range("J68").Value = 49999.2
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 49999.2
$ws.Range("M68").Value = -29189
$ws.Range("N68").Value = -51621.2

$ws.Range("H71").Value = 46666
$ws.Range("I71").Value = 30000
$ws.Range("J71").Value = 49999.2
$ws.Range("K71").Value = 90000
$ws.Range("L71").Value = 149997.6
$ws.Range("M71").Value = -85944
$ws.Range("N71").Value = -158109.6

$ws.Range("H74").Value = 4738.3335
$ws.Range("I74").Value = 4522.857
$ws.Range("K74").Value = 4522.857
$ws.Range("M74").Value = -3648.857

$ws.Range("H77").Value = 4738.3335
$ws.Range("I77").Value = 4522.857
$ws.Range("K77").Value = 22614.285
$ws.Range("M77").Value = -18246.285

$ws.Range("H111").Value = 62055
$ws.Range("J111").Value = 62055
$ws.Range("L111").Value = 62055
$ws.Range("N111").Value = -70235

$ws.Range("H139").Value = 79750
$ws.Range("J139").Value = 79750
$ws.Range("L139").Value = 79750
$ws.Range("N139").Value = -90030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3059.5757
$ws.Range("I3").Value = 3210.1924
$ws.Range("K3").Value = 3210.1924
$ws.Range("M3").Value = -3096.1924

$ws.Range("H80").Value = 472.2857
$ws.Range("I80").Value = 480.33334
$ws.Range("K80").Value = 480.33334
$ws.Range("M80").Value = 517.66666

$ws.Range("H83").Value = 472.2857
$ws.Range("I83").Value = 480.33334
$ws.Range("K83").Value = 2401.6667
$ws.Range("M83").Value = 2590.3333

$ws.Range("H99").Value = 1658.6897
$ws.Range("I99").Value = 1737.5834
$ws.Range("J99").Value = 1280
$ws.Range("K99").Value = 1737.5834
$ws.Range("L99").Value = 1280
$ws.Range("M99").Value = -239.5834
$ws.Range("N99").Value = -4276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 504.29413
$ws.Range("I22").Value = 473.375
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 473.375
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -123.375
$ws.Range("N22").Value = -1699

$ws.Range("H31").Value = 4146.9375
$ws.Range("I31").Value = 2208.25
$ws.Range("J31").Value = 5116.2812
$ws.Range("K31").Value = 2208.25
$ws.Range("L31").Value = 5116.2812
$ws.Range("M31").Value = -1913.25
$ws.Range("N31").Value = -5706.2812

$ws.Range("H34").Value = 4146.9375
$ws.Range("I34").Value = 2208.25
$ws.Range("J34").Value = 5116.2812
$ws.Range("K34").Value = 2208.25
$ws.Range("L34").Value = 5116.2812
$ws.Range("M34").Value = -2006.25
$ws.Range("N34").Value = -5520.2812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3819.8635
$ws.Range("J68").Value = 5670.8335
$ws.Range("L68").Value = 17012.5005
$ws.Range("N68").Value = -18634.5005

$ws.Range("H71").Value = 3819.8635
$ws.Range("J71").Value = 5670.8335
$ws.Range("L71").Value = 51037.5015
$ws.Range("N71").Value = -59149.5015

$ws.Range("H75").Value = 3776.7144
$ws.Range("I75").Value = 1417.4
$ws.Range("J75").Value = 5087.4443
$ws.Range("K75").Value = 4252.200000000001
$ws.Range("L75").Value = 15262.3329
$ws.Range("M75").Value = -3254.200000000001
$ws.Range("N75").Value = -17258.3329

$ws.Range("H78").Value = 3776.7144
$ws.Range("I78").Value = 1417.4
$ws.Range("J78").Value = 5087.4443
$ws.Range("K78").Value = 12756.6
$ws.Range("L78").Value = 45786.9987
$ws.Range("M78").Value = -7764.6
$ws.Range("N78").Value = -55770.9987

$ws.Range("H107").Value = 3747.6
$ws.Range("I107").Value = 372
$ws.Range("J107").Value = 5998
$ws.Range("K107").Value = 1116
$ws.Range("L107").Value = 17994
$ws.Range("M107").Value = 804
$ws.Range("N107").Value = -21834

$ws.Range("H113").Value = 852.0645
$ws.Range("I113").Value = 809.9
$ws.Range("J113").Value = 872.1429000000001
$ws.Range("K113").Value = 2429.7
$ws.Range("L113").Value = 2616.4287
$ws.Range("M113").Value = -259.6999999999998
$ws.Range("N113").Value = -6956.4287

$ws.Range("H119").Value = 14000
$ws.Range("I119").Value = 14000
$ws.Range("K119").Value = 42000
$ws.Range("M119").Value = -37162

$ws.Range("H122").Value = 759.3
$ws.Range("J122").Value = 832.4091
$ws.Range("L122").Value = 7491.6819
$ws.Range("N122").Value = -12391.6819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2627.7693
$ws.Range("I22").Value = 1719.7142
$ws.Range("J22").Value = 3687.1667
$ws.Range("K22").Value = 1719.7142
$ws.Range("L22").Value = 3687.1667
$ws.Range("M22").Value = -1424.7142
$ws.Range("N22").Value = -4277.1667

$ws.Range("H27").Value = 2627.7693
$ws.Range("I27").Value = 1719.7142
$ws.Range("J27").Value = 3687.1667
$ws.Range("K27").Value = 1719.7142
$ws.Range("L27").Value = 3687.1667
$ws.Range("M27").Value = -1612.7142
$ws.Range("N27").Value = -3901.1667

$ws.Range("H39").Value = 5029.5
$ws.Range("I39").Value = 5029.5
$ws.Range("K39").Value = 5029.5
$ws.Range("M39").Value = -4569.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7482.6
$ws.Range("I126").Value = 7482.6
$ws.Range("K126").Value = 22447.8
$ws.Range("M126").Value = -19977.8
